# Slide 60: "Maintaining Context During Parsing"
# The body placeholder's 3rd paragraph currently contains three separate
# runs ("An exit statement has meaning only when nested inside " / "a loop, " /
# "and code generation for an exit statement requires knowledge of which
# loop encloses it.") that together read as one sentence. Merge them into a
# single run (keeping the first run's formatting) without touching the
# paragraph break that follows.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(60)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$para = $tr.Paragraphs(3)

$firstRunText = "An exit statement has meaning only when nested inside "
$firstRunLen = $firstRunText.Length

# Range covering the first run, which already carries the rPr we want to keep.
$firstRun = $tr.Characters($para.Start, $firstRunLen)

# Range covering the remaining two runs' text, excluding the trailing
# paragraph-mark character (Paragraphs().Length includes it).
$restLen = $para.Length - 1 - $firstRunLen
$restRange = $tr.Characters($para.Start + $firstRunLen, $restLen)
$restRange.Delete() | Out-Null

# Re-append the remaining text onto the first run so it keeps that run's
# rPr (lang="en-US" dirty="0") and merges into a single run.
$firstRun.InsertAfter("a loop, and code generation for an exit statement requires knowledge of which loop encloses it.") | Out-Null
